# Applies the scraped-price refresh for Sat Jul 15 13:46:19 UTC 2023.
# Column D ("Price") and E ("Volume(1h)") text is updated per row; row 51
# additionally swaps the coin (Cronos -> SynthetixNetwork) in B/C.
#
# Numeric-looking Price strings (e.g. "250.25") are prefixed with a leading
# apostrophe so Excel stores them as literal text (matching the source sheet,
# where Price is plain text, not a number) instead of silently parsing them
# into floating point numbers and dropping formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.354.53'
$ws.Cells.Item(2, 5).Value = '  -2.73%  '

$ws.Cells.Item(3, 4).Value = '1.937.80'
$ws.Cells.Item(3, 5).Value = '  -2.88%  '

$ws.Cells.Item(4, 5).Value = '  +0.22%  '

$ws.Cells.Item(5, 4).Value = '''250.25'
$ws.Cells.Item(5, 5).Value = '  -1.64%  '

$ws.Cells.Item(6, 4).Value = '''0.7245'
$ws.Cells.Item(6, 5).Value = '  -6.74%  '

$ws.Cells.Item(7, 5).Value = '  +0.18%  '

$ws.Cells.Item(8, 4).Value = '''0.3344'
$ws.Cells.Item(8, 5).Value = '  -3.74%  '

$ws.Cells.Item(9, 4).Value = '''28.38'
$ws.Cells.Item(9, 5).Value = '  +2.25%  '

$ws.Cells.Item(10, 4).Value = '''0.07236'
$ws.Cells.Item(10, 5).Value = '  +2.73%  '

$ws.Cells.Item(11, 4).Value = '''0.8126'
$ws.Cells.Item(11, 5).Value = '  -3.65%  '

$ws.Cells.Item(12, 4).Value = '''0.08098'
$ws.Cells.Item(12, 5).Value = '  -1.05%  '

$ws.Cells.Item(13, 4).Value = '1.936.06'
$ws.Cells.Item(13, 5).Value = '  -2.93%  '

$ws.Cells.Item(14, 4).Value = '''5.476'
$ws.Cells.Item(14, 5).Value = '  -2.61%  '

$ws.Cells.Item(15, 4).Value = '''94.58'
$ws.Cells.Item(15, 5).Value = '  -6.08%  '

$ws.Cells.Item(16, 5).Value = '  -1.57%  '

$ws.Cells.Item(17, 4).Value = '30.375.74'
$ws.Cells.Item(17, 5).Value = '  -2.69%  '

$ws.Cells.Item(18, 4).Value = '''0.000008264'
$ws.Cells.Item(18, 5).Value = '  +3.13%  '

$ws.Cells.Item(19, 4).Value = '''249.70'
$ws.Cells.Item(19, 5).Value = '  -8.30%  '

$ws.Cells.Item(20, 4).Value = '''5.933'
$ws.Cells.Item(20, 5).Value = '  +1.03%  '

$ws.Cells.Item(21, 4).Value = '2.195.16'
$ws.Cells.Item(21, 5).Value = '  -2.57%  '

$ws.Cells.Item(22, 4).Value = '''1.002'
$ws.Cells.Item(22, 5).Value = '  +0.07%  '

$ws.Cells.Item(23, 5).Value = '  +0.22%  '

$ws.Cells.Item(24, 4).Value = '''6.948'
$ws.Cells.Item(24, 5).Value = '  -1.85%  '

$ws.Cells.Item(25, 4).Value = '''9.767'
$ws.Cells.Item(25, 5).Value = '  -2.20%  '

$ws.Cells.Item(26, 4).Value = '''163.54'
$ws.Cells.Item(26, 5).Value = '  -1.23%  '

$ws.Cells.Item(27, 4).Value = '''2.396'
$ws.Cells.Item(27, 5).Value = '  +0.79%  '

$ws.Cells.Item(28, 4).Value = '''19.31'
$ws.Cells.Item(28, 5).Value = '  -2.89%  '

$ws.Cells.Item(29, 4).Value = '''0.1332'
$ws.Cells.Item(29, 5).Value = '  -7.46%  '

$ws.Cells.Item(30, 4).Value = '''1.571'
$ws.Cells.Item(30, 5).Value = '  -1.67%  '

$ws.Cells.Item(31, 4).Value = '''1.347'
$ws.Cells.Item(31, 5).Value = '  -1.52%  '

$ws.Cells.Item(32, 4).Value = '''4.453'
$ws.Cells.Item(32, 5).Value = '  -2.98%  '

$ws.Cells.Item(33, 4).Value = '''4.204'
$ws.Cells.Item(33, 5).Value = '  -4.91%  '

$ws.Cells.Item(34, 4).Value = '''0.05211'
$ws.Cells.Item(34, 5).Value = '  -0.61%  '

$ws.Cells.Item(35, 5).Value = '  +5.90%  '

$ws.Cells.Item(36, 4).Value = '''0.7531'
$ws.Cells.Item(36, 5).Value = '  -4.35%  '

$ws.Cells.Item(37, 4).Value = '''2.752'
$ws.Cells.Item(37, 5).Value = '  -0.52%  '

$ws.Cells.Item(38, 4).Value = '''0.01983'
$ws.Cells.Item(38, 5).Value = '  -0.86%  '

$ws.Cells.Item(39, 4).Value = '''2.843'
$ws.Cells.Item(39, 5).Value = '  -2.48%  '

$ws.Cells.Item(40, 4).Value = '''80.68'
$ws.Cells.Item(40, 5).Value = '  +1.12%  '

$ws.Cells.Item(41, 4).Value = '''6.489'
$ws.Cells.Item(41, 5).Value = '  -3.03%  '

$ws.Cells.Item(42, 4).Value = '''0.4550'

$ws.Cells.Item(43, 4).Value = '''2.041'
$ws.Cells.Item(43, 5).Value = '  -2.73%  '

$ws.Cells.Item(44, 4).Value = '''0.8492'
$ws.Cells.Item(44, 5).Value = '  -0.51%  '

$ws.Cells.Item(45, 5).Value = '  +0.14%  '

$ws.Cells.Item(46, 4).Value = '''102.29'
$ws.Cells.Item(46, 5).Value = '  -2.16%  '

$ws.Cells.Item(47, 4).Value = '''9.803'
$ws.Cells.Item(47, 5).Value = '  -1.59%  '

$ws.Cells.Item(48, 4).Value = '''7.455'
$ws.Cells.Item(48, 5).Value = '  -2.59%  '

$ws.Cells.Item(49, 4).Value = '''36.93'
$ws.Cells.Item(49, 5).Value = '  -0.65%  '

$ws.Cells.Item(50, 4).Value = '''0.4200'
$ws.Cells.Item(50, 5).Value = '  -2.41%  '

$ws.Cells.Item(51, 2).Value = 'SynthetixNetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Cells.Item(51, 4).Value = '''2.871'
$ws.Cells.Item(51, 5).Value = '  +6.38%  '
